$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 2 'Bitcoin'
Set-TextValue 2 3 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextValue 2 4 '25.977.36'
Set-TextValue 2 5 '  +0.66%  '

Set-TextValue 3 2 'Ethereum'
Set-TextValue 3 3 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextValue 3 4 '1.747.77'
Set-TextValue 3 5 '  -0.30%  '

Set-TextValue 4 2 'TetherUSD'
Set-TextValue 4 3 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue 4 4 '1.001'
Set-TextValue 4 5 '  +0.16%  '

Set-TextValue 5 2 'BNB'
Set-TextValue 5 3 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue 5 4 '234.22'
Set-TextValue 5 5 '  -1.19%  '

Set-TextValue 6 2 'USDC'
Set-TextValue 6 3 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue 6 4 '1.000'
Set-TextValue 6 5 '  +0.10%  '

Set-TextValue 7 2 'XRP'
Set-TextValue 7 3 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 7 4 '0.5200'
Set-TextValue 7 5 '  +2.28%  '

Set-TextValue 8 2 'Cardano'
Set-TextValue 8 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 8 4 '0.2823'
Set-TextValue 8 5 '  +4.30%  '

Set-TextValue 9 2 'OKB'
Set-TextValue 9 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 9 4 '39.54'
Set-TextValue 9 5 '  -3.85%  '

Set-TextValue 10 2 'Dogecoin'
Set-TextValue 10 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 10 4 '0.06138'
Set-TextValue 10 5 '  -1.14%  '

Set-TextValue 11 2 'WrappedEther'
Set-TextValue 11 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 11 4 '1.758.03'
Set-TextValue 11 5 '  +1.19%  '

Set-TextValue 12 2 'TRON'
Set-TextValue 12 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 12 4 '0.07022'
Set-TextValue 12 5 '  +1.55%  '

Set-TextValue 13 2 'Solana'
Set-TextValue 13 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 13 4 '15.46'
Set-TextValue 13 5 '  -0.91%  '

Set-TextValue 14 2 'Polygon'
Set-TextValue 14 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 14 4 '0.6450'
Set-TextValue 14 5 '  +2.72%  '

Set-TextValue 15 2 'Polkadot'
Set-TextValue 15 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 15 4 '4.538'
Set-TextValue 15 5 '  +0.92%  '

Set-TextValue 16 2 'Litecoin'
Set-TextValue 16 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 16 4 '77.56'
Set-TextValue 16 5 '  -1.45%  '

Set-TextValue 17 2 'BinanceUSD'
Set-TextValue 17 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 17 4 '1.001'
Set-TextValue 17 5 '  +0.11%  '

Set-TextValue 18 2 'Dai'
Set-TextValue 18 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 18 4 '1.000'
Set-TextValue 18 5 '  +0.05%  '

Set-TextValue 19 2 'WrappedBTC'
Set-TextValue 19 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 19 4 '25.993.37'
Set-TextValue 19 5 '  +0.67%  '

Set-TextValue 20 2 'Avalanche'
Set-TextValue 20 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 20 4 '11.50'
Set-TextValue 20 5 '  -1.80%  '

Set-TextValue 21 2 'ShibaInu'
Set-TextValue 21 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 21 4 '0.000006624'
Set-TextValue 21 5 '  -1.73%  '

Set-TextValue 22 2 'WrappedliquidstakedEther2.0'
Set-TextValue 22 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 22 4 '1.980.54'
Set-TextValue 22 5 '  +0.47%  '

Set-TextValue 23 2 'Uniswap'
Set-TextValue 23 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 23 4 '4.161'
Set-TextValue 23 5 '  +2.23%  '

Set-TextValue 24 2 'Cosmos'
Set-TextValue 24 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 24 4 '8.660'
Set-TextValue 24 5 '  +4.69%  '

Set-TextValue 25 2 'Chainlink'
Set-TextValue 25 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 25 4 '5.155'
Set-TextValue 25 5 '  -0.55%  '

Set-TextValue 26 2 'Monero'
Set-TextValue 26 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 26 4 '139.37'
Set-TextValue 26 5 '  +1.80%  '

Set-TextValue 27 2 'Toncoin'
Set-TextValue 27 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 27 4 '1.509'
Set-TextValue 27 5 '  +3.28%  '

Set-TextValue 28 2 'LidoDAOToken'
Set-TextValue 28 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 28 4 '1.834'
Set-TextValue 28 5 '  +2.04%  '

Set-TextValue 29 2 'EthereumClassic'
Set-TextValue 29 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 29 4 '15.10'
Set-TextValue 29 5 '  -1.60%  '

Set-TextValue 30 2 'BitcoinCash'
Set-TextValue 30 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 30 4 '102.87'
Set-TextValue 30 5 '  +0.15%  '

Set-TextValue 31 2 'Stellar'
Set-TextValue 31 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 31 4 '0.08299'
Set-TextValue 31 5 '  +0.42%  '

Set-TextValue 32 2 'InternetComputer(DFINITY)'
Set-TextValue 32 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 32 4 '3.671'
Set-TextValue 32 5 '  -1.75%  '

Set-TextValue 33 2 'Filecoin'
Set-TextValue 33 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 33 4 '3.445'
Set-TextValue 33 5 '  -0.03%  '

Set-TextValue 34 2 'Hedera'
Set-TextValue 34 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 34 4 '0.04480'
Set-TextValue 34 5 '  +1.53%  '

Set-TextValue 35 2 'HuobiToken'
Set-TextValue 35 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 35 4 '2.616'
Set-TextValue 35 5 '  -1.05%  '

Set-TextValue 36 2 'ARBITRUM'
Set-TextValue 36 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 36 4 '0.9892'
Set-TextValue 36 5 '  -1.89%  '

Set-TextValue 37 2 'ImmutableX'
Set-TextValue 37 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 37 4 '0.6158'
Set-TextValue 37 5 '  +1.54%  '

Set-TextValue 38 2 'MXToken'
Set-TextValue 38 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 38 4 '2.684'
Set-TextValue 38 5 '  -0.55%  '

Set-TextValue 39 2 'VeChain'
Set-TextValue 39 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 39 4 '0.01590'
Set-TextValue 39 5 '  +1.91%  '

Set-TextValue 40 2 'RenderToken'
Set-TextValue 40 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 40 4 '1.936'
Set-TextValue 40 5 '  -1.70%  '

Set-TextValue 41 2 'PaxDollar'
Set-TextValue 41 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 41 4 '0.9997'
Set-TextValue 41 5 '  -0.03%  '

Set-TextValue 42 2 'PaxosStandard'
Set-TextValue 42 3 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
Set-TextValue 42 4 '1.001'
Set-TextValue 42 5 '  +0.08%  '

Set-TextValue 43 2 'Quant'
Set-TextValue 43 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 43 4 '100.69'
Set-TextValue 43 5 '  -1.47%  '

Set-TextValue 44 2 'TheSandbox'
Set-TextValue 44 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 44 4 '0.3871'
Set-TextValue 44 5 '  +0.11%  '

Set-TextValue 45 2 'FraxShare'
Set-TextValue 45 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 45 4 '5.089'
Set-TextValue 45 5 '  +4.43%  '

Set-TextValue 46 2 'TrustWalletToken'
Set-TextValue 46 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 46 4 '0.7365'
Set-TextValue 46 5 '  -1.91%  '

Set-TextValue 47 2 'Cronos'
Set-TextValue 47 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 47 4 '0.05471'
Set-TextValue 47 5 '  -0.71%  '

Set-TextValue 48 2 'Aptos'
Set-TextValue 48 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 48 4 '6.312'
Set-TextValue 48 5 '  +5.69%  '

Set-TextValue 49 2 'Algorand'
Set-TextValue 49 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 49 4 '0.1127'
Set-TextValue 49 5 '  +2.83%  '

Set-TextValue 50 2 'Aave'
Set-TextValue 50 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 50 4 '53.04'
Set-TextValue 50 5 '  +0.44%  '

Set-TextValue 51 2 'Elrond'
Set-TextValue 51 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 51 4 '30.04'
Set-TextValue 51 5 '  -0.85%  '

